$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 115
$ws.Range("H115").Value = 3436.8
$ws.Range("I115").Value = 3728
$ws.Range("K115").Value = 11184
$ws.Range("M115").Value = -9617
# Row 132
$ws.Range("H132").Value = 4264.8
$ws.Range("I132").Value = 4389.5386
$ws.Range("J132").Value = 3904.4443
$ws.Range("K132").Value = 13168.6158
$ws.Range("L132").Value = 11713.3329
$ws.Range("M132").Value = -10638.6158
$ws.Range("N132").Value = -16773.3329
# Row 138
$ws.Range("H138").Value = 2550.1777
$ws.Range("I138").Value = 3099.7693
$ws.Range("J138").Value = 2457.3896
$ws.Range("K138").Value = 9299.3079
$ws.Range("L138").Value = 7372.168799999999
$ws.Range("M138").Value = -4159.3079
$ws.Range("N138").Value = -17652.1688

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2214.8667
$ws.Range("I2").Value = 2338.4546
$ws.Range("J2").Value = 1875
$ws.Range("K2").Value = 2338.4546
$ws.Range("L2").Value = 1875
$ws.Range("M2").Value = -2225.4546
$ws.Range("N2").Value = -2101
# Row 32
$ws.Range("H32").Value = 420244.72
$ws.Range("I32").Value = 465373.94
$ws.Range("K32").Value = 465373.94
$ws.Range("M32").Value = -465086.94
# Row 45
$ws.Range("H45").Value = 2386.4443
$ws.Range("I45").Value = 2264.182
$ws.Range("J45").Value = 2470.5
$ws.Range("K45").Value = 2264.182
$ws.Range("L45").Value = 2470.5
$ws.Range("M45").Value = -1887.182
$ws.Range("N45").Value = -3224.5
# Row 61
$ws.Range("H61").Value = 3097.7666
$ws.Range("I61").Value = 2658.2942
$ws.Range("K61").Value = 2658.2942
$ws.Range("M61").Value = -2446.2942
# Row 68
$ws.Range("H68").Value = 47500
$ws.Range("J68").Value = 47500
$ws.Range("L68").Value = 47500
$ws.Range("N68").Value = -49122
# Row 71
$ws.Range("H71").Value = 47500
$ws.Range("J71").Value = 47500
$ws.Range("L71").Value = 142500
$ws.Range("N71").Value = -150612
# Row 74
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()
# Row 77
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()
# Row 110
$ws.Range("H110").Value = 1323.8462
$ws.Range("I110").Value = 1323.8462
$ws.Range("K110").Value = 1323.8462
$ws.Range("M110").Value = 721.1538
# Row 116
$ws.Range("H116").Value = 2214.8667
$ws.Range("I116").Value = 2338.4546
$ws.Range("J116").Value = 1875
$ws.Range("K116").Value = 2338.4546
$ws.Range("L116").Value = 1875
$ws.Range("M116").Value = -44.45460000000003
$ws.Range("N116").Value = -6463
# Row 136
$ws.Range("H136").Value = 3097.7666
$ws.Range("I136").Value = 2658.2942
$ws.Range("K136").Value = 7974.882599999999
$ws.Range("M136").Value = -5424.882599999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2214.8667
$ws.Range("I3").Value = 2338.4546
$ws.Range("J3").Value = 1875
$ws.Range("K3").Value = 2338.4546
$ws.Range("L3").Value = 1875
$ws.Range("M3").Value = -2224.4546
$ws.Range("N3").Value = -2103
# Row 99
$ws.Range("H99").Value = 913
$ws.Range("I99").Value = 803.8333
$ws.Range("J99").Value = 1032.091
$ws.Range("K99").Value = 803.8333
$ws.Range("L99").Value = 1032.091
$ws.Range("M99").Value = 694.1667
$ws.Range("N99").Value = -4028.091
# Row 105
$ws.Range("H105").Value = 11366843
$ws.Range("I105").Value = 15628085
$ws.Range("K105").Value = 15628085
$ws.Range("M105").Value = -15626338
# Row 107
$ws.Range("H107").Value = 44575.22
$ws.Range("I107").Value = 53674.633
$ws.Range("K107").Value = 53674.633
$ws.Range("M107").Value = -51754.633
# Row 134
$ws.Range("H134").Value = 2579.2942
$ws.Range("J134").Value = 2950.9285
$ws.Range("L134").Value = 8852.7855
$ws.Range("N134").Value = -13922.7855

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5752.163
$ws.Range("I31").Value = 1182.4828
$ws.Range("K31").Value = 1182.4828
$ws.Range("M31").Value = -887.4828
# Row 34
$ws.Range("H34").Value = 5752.163
$ws.Range("I34").Value = 1182.4828
$ws.Range("K34").Value = 1182.4828
$ws.Range("M34").Value = -980.4828
# Row 62
$ws.Range("H62").Value = 4436.364
$ws.Range("J62").Value = 4500
$ws.Range("L62").Value = 4500
$ws.Range("N62").Value = -5748
# Row 65
$ws.Range("H65").Value = 4436.364
$ws.Range("J65").Value = 4500
$ws.Range("L65").Value = 22500
$ws.Range("N65").Value = -28740
# Row 105
$ws.Range("H105").Value = 1981.1428
$ws.Range("J105").Value = 1950
$ws.Range("L105").Value = 1950
$ws.Range("N105").Value = -5444
# Row 132
$ws.Range("H132").Value = 7939132.5
$ws.Range("I132").Value = 2492.2
$ws.Range("J132").Value = 15154260
$ws.Range("K132").Value = 7476.599999999999
$ws.Range("L132").Value = 45462780
$ws.Range("M132").Value = -4946.599999999999
$ws.Range("N132").Value = -45467840

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 200
$ws.Range("I11").Value = 200
$ws.Range("K11").Value = 600
$ws.Range("M11").Value = -460

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 2563.7188
$ws.Range("I132").Value = 1544.3158
$ws.Range("K132").Value = 4632.9474
$ws.Range("M132").Value = -2102.9474

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 43
$ws.Range("H43").Value = 3000
$ws.Range("I43").Value = 3000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 3000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -2851
$ws.Range("N43").ClearContents()
# Row 74
$ws.Range("H74").Value = 10400
$ws.Range("J74").Value = 10600
$ws.Range("L74").Value = 10600
$ws.Range("N74").Value = -12472
# Row 77
$ws.Range("H77").Value = 10400
$ws.Range("J77").Value = 10600
$ws.Range("L77").Value = 31800
$ws.Range("N77").Value = -41160
# Row 132
$ws.Range("H132").Value = 4168925.8
$ws.Range("I132").Value = 2308.4092
$ws.Range("J132").Value = 9261458
$ws.Range("K132").Value = 6925.2276
$ws.Range("L132").Value = 27784374
$ws.Range("M132").Value = -4395.2276
$ws.Range("N132").Value = -27789434
